$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The body of the table (rows 10-25) gets reshuffled: several "label" rows
# absorb the value that used to sit two rows below them, a few rows lose
# their B/C value altogether, and the last two rows collapse away once
# their content has been relocated. We copy cell-to-cell (instead of typing
# literal values) so Excel keeps reusing the existing shared-string entries
# and cell styles exactly as a manual cut/paste edit would, and so that
# text which merely *looks* like a date ("01/01/2012") is not reinterpreted
# as a date serial number.

# Snapshot the two B/C values that would otherwise be overwritten before we
# need to reuse them further down (B13/C13 feed both B10/C10 and B15/C15).
$ws.Range("B13").Copy()
$ws.Range("B10").PasteSpecial(-4104)
$ws.Range("C13").Copy()
$ws.Range("C10").PasteSpecial(-4104)

$ws.Range("B13").Copy()
$ws.Range("B15").PasteSpecial(-4104)
$ws.Range("C13").Copy()
$ws.Range("C15").PasteSpecial(-4104)

# Row 13 becomes the "Programa resumido:" row, reusing the "01/01/2012" text.
$ws.Range("A15").Copy()
$ws.Range("A13").PasteSpecial(-4104)
$ws.Range("B8").Copy()
$ws.Range("B13").PasteSpecial(-4104)
$ws.Range("C8").Copy()
$ws.Range("C13").PasteSpecial(-4104)

# Row 14 becomes "Short syllabus:" only.
$ws.Range("A16").Copy()
$ws.Range("A14").PasteSpecial(-4104)

# Row 15's label becomes "Programa:" (value already placed above).
$ws.Range("A17").Copy()
$ws.Range("A15").PasteSpecial(-4104)

# Row 16 becomes "Syllabus:" only.
$ws.Range("A18").Copy()
$ws.Range("A16").PasteSpecial(-4104)

# Row 17 becomes "Avaliação:" only.
$ws.Range("A19").Copy()
$ws.Range("A17").PasteSpecial(-4104)

# Row 18 becomes the "Método:" row.
$ws.Range("A20").Copy()
$ws.Range("A18").PasteSpecial(-4104)
$ws.Range("B14").Copy()
$ws.Range("B18").PasteSpecial(-4104)
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4104)

# Row 19 becomes the "Critério:" row.
$ws.Range("A21").Copy()
$ws.Range("A19").PasteSpecial(-4104)
$ws.Range("B20").Copy()
$ws.Range("B19").PasteSpecial(-4104)
$ws.Range("C20").Copy()
$ws.Range("C19").PasteSpecial(-4104)

# Row 20 becomes the "Norma de recuperação:" row.
$ws.Range("A22").Copy()
$ws.Range("A20").PasteSpecial(-4104)
$ws.Range("B21").Copy()
$ws.Range("B20").PasteSpecial(-4104)
$ws.Range("C21").Copy()
$ws.Range("C20").PasteSpecial(-4104)

# Row 21 becomes the "Bibliografia:" row.
$ws.Range("A23").Copy()
$ws.Range("A21").PasteSpecial(-4104)
$ws.Range("B22").Copy()
$ws.Range("B21").PasteSpecial(-4104)
$ws.Range("C22").Copy()
$ws.Range("C21").PasteSpecial(-4104)

# Row 22 becomes "Requisitos:" only.
$ws.Range("A24").Copy()
$ws.Range("A22").PasteSpecial(-4104)

# Row 23 keeps only the B/C requirement text (no A label).
$ws.Range("B25").Copy()
$ws.Range("B23").PasteSpecial(-4104)
$ws.Range("C25").Copy()
$ws.Range("C23").PasteSpecial(-4104)

# Clear the cells that must now be empty on rows 14, 17, 22 and 23.
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Range("A23").Clear()

# Row heights for the reshuffled rows. Rows 17 and 22 go back to the
# (non-custom) default height, so use AutoFit to drop the customHeight flag
# instead of pinning a literal RowHeight.
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30

# Rows 24 and 25 are now obsolete (their content already relocated above).
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(24).Delete()
